$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.99850114534168311
$ws.Range("A2").Value = 0.99928035316489638
$ws.Range("A3").Value = 0.76556987424224476
$ws.Range("A4").Value = 0.63981047638672361
$ws.Range("A5").Value = 0.49438961894859962
$ws.Range("A6").Value = 0.55442968963806993
$ws.Range("A7").Value = 0.56094725214672314
$ws.Range("A8").Value = 0.78494766785344106
$ws.Range("A9").Value = 0.93389049251852485
$ws.Range("A10").Value = 0.6305070102378475
$ws.Range("A11").Value = 0.67810123712965431
$ws.Range("A12").Value = 0.46067451238762586
$ws.Range("A13").Value = 0.46612228605497885
$ws.Range("A14").Value = 0.53852243141626877
$ws.Range("A15").Value = 0.20299709768542135
$ws.Range("A16").Value = 0.69281855351465949
$ws.Range("A17").Value = 0.2647152327009385
$ws.Range("A18").Value = 0.39889824572679639
$ws.Range("A19").Value = 0.37593146482017287
$ws.Range("A20").Value = 0.3624486942936726
$ws.Range("A21").Value = 0.37909743397727136
$ws.Range("A22").Value = 0.11542716448814157
$ws.Range("A23").Value = 0.23332720634932536
$ws.Range("A24").Value = 0.32914002361581568
$ws.Range("A25").Value = 0.23307982285222506
$ws.Range("A26").Value = 0.47586956981582873
$ws.Range("A27").Value = 0.73495681001872915
$ws.Range("A28").Value = 0.53845201758012229
$ws.Range("A29").Value = 0.42604259268477357
$ws.Range("A30").Value = 0.37705119766408929
$ws.Range("A31").Value = 0.26769039847228138
$ws.Range("A32").Value = 0.3958494147432568
$ws.Range("A33").Value = 0.43446264323303208
$ws.Range("A34").Value = 0.45029230464351055
$ws.Range("A35").Value = 0.61132604257110346
$ws.Range("A36").Value = 0.60626961110925381
$ws.Range("A37").Value = 0.61845464606076372
$ws.Range("A38").Value = 0.65713961382304109
$ws.Range("A39").Value = 0.74138401486627503
$ws.Range("A40").Value = 0.70214817024409393
$ws.Range("A41").Value = 0.52021201184863541
$ws.Range("A42").Value = 0.60535045189067183
$ws.Range("A43").Value = 0.48821164007359519
$ws.Range("A44").Value = 0.47977011115672097
$ws.Range("A45").Value = 0.4749637132450743
$ws.Range("A46").Value = 0.37064462817964716
$ws.Range("A47").Value = 0.37138974227303101
$ws.Range("A48").Value = 0.39001148449285694
$ws.Range("A49").Value = 0.50296804482521407
$ws.Range("A50").Value = 0.40323239081410978
$ws.Range("A51").Value = 0.45072211248275074
$ws.Range("A52").Value = 0.43687182600181146
$ws.Range("A53").Value = 0.53047946438166649
$ws.Range("A54").Value = 0.53381936152368203
$ws.Range("A55").Value = 0.55265692800978516
$ws.Range("A56").Value = 0.60329978759780123
$ws.Range("A57").Value = 0.63389499940373761
$ws.Range("A58").Value = 0.52086263021228885
$ws.Range("A59").Value = 0.48096398230913701
